# Applies updated "dSF" (column F) values for eovaldi_nathan.xlsx
# Commit message: "repull data, push all data, mean calculation"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -1
    3  = -1
    4  = 6
    5  = 1
    6  = 2
    7  = 3
    8  = -6
    9  = 2
    10 = -1
    13 = 0
    14 = 2
    15 = -2
    16 = 13
    17 = -3
    18 = 0
    19 = 3
    21 = -2
    22 = 12
    23 = 4
    24 = -2
    25 = 1
    26 = -2
    28 = 2
    29 = 2
    30 = -3
    32 = 1
    33 = 5
    34 = -5
    35 = 7
    36 = 1
    37 = 8
    39 = 2
    40 = -1
    41 = -4
    42 = 1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
